# Auto-generated edit script applying the Pandaemonium_Profits.xlsx diff
# Updates currentAveragePrice* / Leve* columns (H:N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N62").Value2 = -8676.5713  # was -7935.5
$ws.Range("L62").Value2 = 7428.5713  # was 6687.5
$ws.Range("J62").Value2 = 7428.5713  # was 6687.5
$ws.Range("M62").Value2 = -1759.75  # was -1894.6365
$ws.Range("H62").Value2 = 4242.3687  # was 4273.9473
$ws.Range("K62").Value2 = 2383.75  # was 2518.6365
$ws.Range("I62").Value2 = 2383.75  # was 2518.6365
$ws.Range("K65").Value2 = 11918.75  # was 12593.1825
$ws.Range("L65").Value2 = 37142.85649999999  # was 33437.5
$ws.Range("I65").Value2 = 2383.75  # was 2518.6365
$ws.Range("M65").Value2 = -8798.75  # was -9473.182500000001
$ws.Range("H65").Value2 = 4242.3687  # was 4273.9473
$ws.Range("N65").Value2 = -43382.85649999999  # was -39677.5
$ws.Range("J65").Value2 = 7428.5713  # was 6687.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value2 = 5670.3335  # was 3768.5
$ws.Range("H2").Value2 = 5670.3335  # was 2840.9092
$ws.Range("J2").Value2 = 0  # was 1727.8
$ws.Range("N2").ClearContents()  # was -1953.8
$ws.Range("K2").Value2 = 5670.3335  # was 3768.5
$ws.Range("L2").Value2 = 0  # was 1727.8
$ws.Range("M2").Value2 = -5557.3335  # was -3655.5
$ws.Range("H97").Value2 = 704.23334  # was 881.38464
$ws.Range("K97").Value2 = 604.46155  # was 700.7619
$ws.Range("J97").Value2 = 1352.75  # was 1640
$ws.Range("I97").Value2 = 604.46155  # was 700.7619
$ws.Range("N97").Value2 = -2344.75  # was -2632
$ws.Range("L97").Value2 = 1352.75  # was 1640
$ws.Range("M97").Value2 = -108.46155  # was -204.7619
$ws.Range("I110").Value2 = 1466  # was 1561.4286
$ws.Range("K110").Value2 = 1466  # was 1561.4286
$ws.Range("M110").Value2 = 579  # was 483.5714
$ws.Range("L110").Value2 = 1996  # was 2120
$ws.Range("J110").Value2 = 1996  # was 2120
$ws.Range("N110").Value2 = -6086  # was -6210
$ws.Range("H110").Value2 = 1669.8462  # was 1764.5454
$ws.Range("N116").ClearContents()  # was -6315.8
$ws.Range("L116").Value2 = 0  # was 1727.8
$ws.Range("M116").Value2 = -3376.3335  # was -1474.5
$ws.Range("H116").Value2 = 5670.3335  # was 2840.9092
$ws.Range("J116").Value2 = 0  # was 1727.8
$ws.Range("K116").Value2 = 5670.3335  # was 3768.5
$ws.Range("I116").Value2 = 5670.3335  # was 3768.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value2 = 0  # was 1727.8
$ws.Range("H3").Value2 = 5670.3335  # was 2840.9092
$ws.Range("I3").Value2 = 5670.3335  # was 3768.5
$ws.Range("N3").ClearContents()  # was -1955.8
$ws.Range("L3").Value2 = 0  # was 1727.8
$ws.Range("M3").Value2 = -5556.3335  # was -3654.5
$ws.Range("K3").Value2 = 5670.3335  # was 3768.5
$ws.Range("H20").Value2 = 7491.5  # was 7672
$ws.Range("M20").Value2 = -7244.5  # was -7425
$ws.Range("K20").Value2 = 7491.5  # was 7672
$ws.Range("I20").Value2 = 7491.5  # was 7672
$ws.Range("H80").Value2 = 247.52942  # was 186.26315
$ws.Range("I80").Value2 = 103.666664  # was 147.5
$ws.Range("M80").Value2 = 894.333336  # was 850.5
$ws.Range("J80").Value2 = 278.35715  # was 190.82353
$ws.Range("L80").Value2 = 278.35715  # was 190.82353
$ws.Range("N80").Value2 = -2274.35715  # was -2186.82353
$ws.Range("K80").Value2 = 103.666664  # was 147.5
$ws.Range("L83").Value2 = 1391.78575  # was 954.11765
$ws.Range("M83").Value2 = 4473.66668  # was 4254.5
$ws.Range("H83").Value2 = 247.52942  # was 186.26315
$ws.Range("I83").Value2 = 103.666664  # was 147.5
$ws.Range("J83").Value2 = 278.35715  # was 190.82353
$ws.Range("N83").Value2 = -11375.78575  # was -10938.11765
$ws.Range("K83").Value2 = 518.33332  # was 737.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K31").Value2 = 1575.5  # was 1538.25
$ws.Range("L31").Value2 = 2959.1428  # was 3000
$ws.Range("J31").Value2 = 2959.1428  # was 3000
$ws.Range("N31").Value2 = -3549.1428  # was -3590
$ws.Range("H31").Value2 = 2036.7142  # was 1830.6
$ws.Range("I31").Value2 = 1575.5  # was 1538.25
$ws.Range("M31").Value2 = -1280.5  # was -1243.25
$ws.Range("K34").Value2 = 1575.5  # was 1538.25
$ws.Range("H34").Value2 = 2036.7142  # was 1830.6
$ws.Range("L34").Value2 = 2959.1428  # was 3000
$ws.Range("I34").Value2 = 1575.5  # was 1538.25
$ws.Range("J34").Value2 = 2959.1428  # was 3000
$ws.Range("M34").Value2 = -1373.5  # was -1336.25
$ws.Range("N34").Value2 = -3363.1428  # was -3404
$ws.Range("H82").Value2 = 34590.5  # was 32590.5
$ws.Range("N82").Value2 = -35842.668  # was -33312.5
$ws.Range("M82").Value2 = -32639  # was None
$ws.Range("J82").Value2 = 35120.668  # was 32590.5
$ws.Range("K82").Value2 = 33000  # was 0
$ws.Range("I82").Value2 = 33000  # was 0
$ws.Range("L82").Value2 = 35120.668  # was 32590.5
$ws.Range("H85").Value2 = 34590.5  # was 32590.5
$ws.Range("M85").Value2 = -31752  # was None
$ws.Range("I85").Value2 = 33000  # was 0
$ws.Range("N85").Value2 = -37616.668  # was -35086.5
$ws.Range("L85").Value2 = 35120.668  # was 32590.5
$ws.Range("J85").Value2 = 35120.668  # was 32590.5
$ws.Range("K85").Value2 = 33000  # was 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value2 = 3998  # was 3999.25
$ws.Range("N82").Value2 = -12806  # was -12809.75
$ws.Range("J82").Value2 = 3998  # was 3999.25
$ws.Range("L82").Value2 = 11994  # was 11997.75
$ws.Range("H85").Value2 = 3998  # was 3999.25
$ws.Range("N85").Value2 = -14802  # was -14805.75
$ws.Range("L85").Value2 = 11994  # was 11997.75
$ws.Range("J85").Value2 = 3998  # was 3999.25
$ws.Range("J98").Value2 = 813.4286  # was 918.8
$ws.Range("L98").Value2 = 2440.2858  # was 2756.4
$ws.Range("H98").Value2 = 442.89655  # was 434.96295
$ws.Range("N98").Value2 = -5436.2858  # was -5752.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 0  # was 1090.4
$ws.Range("K97").Value2 = 0  # was 1006.75
$ws.Range("J97").Value2 = 0  # was 1239.1111
$ws.Range("I97").Value2 = 0  # was 1006.75
$ws.Range("N97").ClearContents()  # was -2231.1111
$ws.Range("L97").Value2 = 0  # was 1239.1111
$ws.Range("M97").ClearContents()  # was -510.75
$ws.Range("H102").Value2 = 4089.6765  # was 3755.5789
$ws.Range("I102").Value2 = 3744.88  # was 3354.6553
$ws.Range("K102").Value2 = 3744.88  # was 3354.6553
$ws.Range("M102").Value2 = -2122.88  # was -1732.6553
$ws.Range("H122").Value2 = 7399.8887  # was 13753.25
$ws.Range("I122").Value2 = 11259.8  # was 13753.25
$ws.Range("J122").Value2 = 2575  # was 0
$ws.Range("L122").Value2 = 7725  # was 0
$ws.Range("M122").Value2 = -31329.39999999999  # was -38809.75
$ws.Range("N122").Value2 = -12625  # was None
$ws.Range("K122").Value2 = 33779.39999999999  # was 41259.75
$ws.Range("N125").Value2 = -84720  # was -34920
$ws.Range("J125").Value2 = 79800  # was 30000
$ws.Range("L125").Value2 = 79800  # was 30000
$ws.Range("H125").Value2 = 79800  # was 30000
$ws.Range("M126").Value2 = -3539  # was -3452.800099999999
$ws.Range("K126").Value2 = 6009  # was 5922.800099999999
$ws.Range("N126").Value2 = -16873.3334  # was -15849.0905
$ws.Range("L126").Value2 = 11933.3334  # was 10909.0905
$ws.Range("I126").Value2 = 2003  # was 1974.2667
$ws.Range("H126").Value2 = 3048.4707  # was 2677.4614
$ws.Range("J126").Value2 = 3977.7778  # was 3636.3635
$ws.Range("J127").Value2 = 48900  # was 30000
$ws.Range("N127").Value2 = -58820  # was -39920
$ws.Range("H127").Value2 = 48900  # was 30000
$ws.Range("L127").Value2 = 48900  # was 30000
$ws.Range("L128").Value2 = 79800  # was 0
$ws.Range("H128").Value2 = 79800  # was 0
$ws.Range("J128").Value2 = 79800  # was 0
$ws.Range("N128").Value2 = -89760  # was None
$ws.Range("J129").Value2 = 47449.5  # was 49932.668
$ws.Range("H129").Value2 = 47449.5  # was 49932.668
$ws.Range("N129").Value2 = -57449.5  # was -59932.668
$ws.Range("L129").Value2 = 47449.5  # was 49932.668
$ws.Range("H130").Value2 = 72257.39999999999  # was 79800
$ws.Range("N130").Value2 = -82297.39999999999  # was -89840
$ws.Range("J130").Value2 = 72257.39999999999  # was 79800
$ws.Range("L130").Value2 = 72257.39999999999  # was 79800

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J7").Value2 = 11125  # was 4350.75
$ws.Range("K7").Value2 = 5469.4  # was 2333.8
$ws.Range("L7").Value2 = 11125  # was 4350.75
$ws.Range("M7").Value2 = -5357.4  # was -2221.8
$ws.Range("H7").Value2 = 7085.2856  # was 2910.0715
$ws.Range("N7").Value2 = -11349  # was -4574.75
$ws.Range("I7").Value2 = 5469.4  # was 2333.8
$ws.Range("I40").Value2 = 3517.5833  # was 3662.2307
$ws.Range("H40").Value2 = 3826.7368  # was 4006.8125
$ws.Range("N40").Value2 = -4628.7144  # was -5772
$ws.Range("K40").Value2 = 3517.5833  # was 3662.2307
$ws.Range("M40").Value2 = -3381.5833  # was -3526.2307
$ws.Range("J40").Value2 = 4356.7144  # was 5500
$ws.Range("L40").Value2 = 4356.7144  # was 5500
$ws.Range("H82").Value2 = 1516.5  # was 1620
$ws.Range("M82").Value2 = -888.5  # was -1139
$ws.Range("K82").Value2 = 1249.5  # was 1500
$ws.Range("I82").Value2 = 1249.5  # was 1500
$ws.Range("H85").Value2 = 1516.5  # was 1620
$ws.Range("M85").Value2 = -1.5  # was -252
$ws.Range("I85").Value2 = 1249.5  # was 1500
$ws.Range("K85").Value2 = 1249.5  # was 1500
$ws.Range("L94").Value2 = 24166.666  # was 21875
$ws.Range("H94").Value2 = 24166.666  # was 21875
$ws.Range("J94").Value2 = 24166.666  # was 21875
$ws.Range("N94").Value2 = -25518.666  # was -23227
$ws.Range("M126").Value2 = -13938.2  # was -4531.400000000001
$ws.Range("K126").Value2 = 16408.2  # was 7001.400000000001
$ws.Range("N126").Value2 = -38315  # was -17992.25
$ws.Range("L126").Value2 = 33375  # was 13052.25
$ws.Range("I126").Value2 = 5469.4  # was 2333.8
$ws.Range("H126").Value2 = 7085.2856  # was 2910.0715
$ws.Range("J126").Value2 = 11125  # was 4350.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K107").Value2 = 1776  # was 1986
$ws.Range("I107").Value2 = 592  # was 662
$ws.Range("H107").Value2 = 3232.1538  # was 3264.4614
$ws.Range("M107").Value2 = 144  # was -66
